$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.006.34"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "2.118.35"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "348.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5203"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4456"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09357"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.452"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.883"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.090.47"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06695"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.309"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.008"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").Value = "30.037.49"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").Value = "2.397.07"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.550"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.154"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.780"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.255"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.570"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02606"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06869"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7040"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.336"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2250"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6862"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.353"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.42%  "

$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000360"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.636"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.225"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
